$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J (copy header formatting from H1, then set text)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I and J, rows 2-30
$data = @{
    2  = @(1, 6)
    3  = @(1, 5)
    4  = @(1, 4)
    5  = @(1, 5)
    6  = @(1, 6)
    7  = @(1, 6)
    8  = @(1, 6)
    9  = @(1, 6)
    10 = @(1, 7)
    11 = @(1, 5)
    12 = @(1, 6)
    13 = @(1, 6)
    14 = @(1, 6)
    15 = @(1, 3)
    16 = @(1, 6)
    17 = @(1, 3)
    18 = @(1, 6)
    19 = @(1, 5)
    20 = @(1, 6)
    21 = @(1, 7)
    22 = @(1, 5)
    23 = @(1, 5)
    24 = @(1, 6)
    25 = @(1, 3)
    26 = @(1, 4)
    27 = @(1, 3)
    28 = @(5, 6)
    29 = @(3, 4)
    30 = @(1, 1)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
